$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DValue($sheet, $row, $value, $isNumericLooking) {
    $cell = $sheet.Range("D" + $row)
    if ($isNumericLooking) {
        $cell.Value = "'" + $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}

Set-DValue $ws 2 "63.390.18" $false
$ws.Range("E2").Value = "  -7.33%  "

Set-DValue $ws 3 "3.280.12" $false
$ws.Range("E3").Value = "  -8.38%  "

$ws.Range("E4").Value = "  -0.06%  "

Set-DValue $ws 5 "181.80" $true
$ws.Range("E5").Value = "  -12.27%  "

Set-DValue $ws 6 "516.70" $true
$ws.Range("E6").Value = "  -9.03%  "

Set-DValue $ws 7 "0.595" $true
$ws.Range("E7").Value = "  -2.31%  "

Set-DValue $ws 8 "3.277.51" $false
$ws.Range("E8").Value = "  -8.34%  "

$ws.Range("E9").Value = "  -0.03%  "

Set-DValue $ws 10 "0.618" $true
$ws.Range("E10").Value = "  -8.63%  "

Set-DValue $ws 11 "59.07" $true
$ws.Range("E11").Value = "  -7.13%  "

Set-DValue $ws 12 "0.131" $true
$ws.Range("E12").Value = "  -10.99%  "

Set-DValue $ws 13 "0.0000255" $true
$ws.Range("E13").Value = "  -9.44%  "

Set-DValue $ws 14 "9.11" $true
$ws.Range("E14").Value = "  -10.09%  "

Set-DValue $ws 15 "3.817.46" $false
$ws.Range("E15").Value = "  -8.27%  "

Set-DValue $ws 16 "0.118" $true
$ws.Range("E16").Value = "  -5.51%  "

Set-DValue $ws 17 "3.296.58" $false
$ws.Range("E17").Value = "  -8.43%  "

Set-DValue $ws 18 "17.53" $true
$ws.Range("E18").Value = "  -8.49%  "

Set-DValue $ws 19 "63.346.37" $false
$ws.Range("E19").Value = "  -7.12%  "

Set-DValue $ws 20 "10.95" $true
$ws.Range("E20").Value = "  -10.20%  "

Set-DValue $ws 21 "0.946" $true
$ws.Range("E21").Value = "  -10.87%  "

Set-DValue $ws 22 "371.17" $true
$ws.Range("E22").Value = "  -8.36%  "

Set-DValue $ws 23 "11.20" $true
$ws.Range("E23").Value = "  -9.31%  "

Set-DValue $ws 24 "80.20" $true
$ws.Range("E24").Value = "  -5.31%  "

Set-DValue $ws 25 "3.66" $true
$ws.Range("E25").Value = "  -11.73%  "

Set-DValue $ws 26 "3.88" $true
$ws.Range("E26").Value = "  +1.62%  "

Set-DValue $ws 27 "5.93" $true
$ws.Range("E27").Value = "  -3.41%  "

Set-DValue $ws 28 "2.65" $true
$ws.Range("E28").Value = "  -7.75%  "

Set-DValue $ws 29 "11.37" $true
$ws.Range("E29").Value = "  -8.68%  "

Set-DValue $ws 30 "8.31" $true
$ws.Range("E30").Value = "  -8.50%  "

Set-DValue $ws 31 "647.66" $true
$ws.Range("E31").Value = "  -11.87%  "

Set-DValue $ws 32 "28.44" $true
$ws.Range("E32").Value = "  -9.60%  "

Set-DValue $ws 33 "6.68" $true
$ws.Range("E33").Value = "  -11.48%  "

Set-DValue $ws 34 "11.17" $true
$ws.Range("E34").Value = "  -7.40%  "

Set-DValue $ws 35 "59.57" $true
$ws.Range("E35").Value = "  -6.83%  "

Set-DValue $ws 36 "0.105" $true
$ws.Range("E36").Value = "  -6.56%  "

Set-DValue $ws 37 "0.999" $true
$ws.Range("E37").Value = "  -0.02%  "

Set-DValue $ws 38 "0.387" $true
$ws.Range("E38").Value = "  -8.35%  "

Set-DValue $ws 39 "36.02" $true
$ws.Range("E39").Value = "  -12.79%  "

$ws.Range("E40").Value = "  -0.02%  "

Set-DValue $ws 41 "2.980.15" $false
$ws.Range("E41").Value = "  -5.83%  "

Set-DValue $ws 42 "0.125" $true
$ws.Range("E42").Value = "  -5.80%  "

Set-DValue $ws 43 "0.0₃0650" $false
$ws.Range("E43").Value = "  -11.66%  "

Set-DValue $ws 44 "2.68" $true
$ws.Range("E44").Value = "  -17.24%  "

Set-DValue $ws 45 "2.39" $true
$ws.Range("E45").Value = "  -7.29%  "

$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-DValue $ws 46 "2.59" $true
$ws.Range("E46").Value = "  -5.76%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-DValue $ws 47 "0.0388" $true
$ws.Range("E47").Value = "  -5.66%  "

Set-DValue $ws 48 "2.79" $true
$ws.Range("E48").Value = "  +2.92%  "

Set-DValue $ws 49 "0.125" $true
$ws.Range("E49").Value = "  -4.01%  "

Set-DValue $ws 50 "2.95" $true
$ws.Range("E50").Value = "  -4.44%  "

$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-DValue $ws 51 "131.31" $true
$ws.Range("E51").Value = "  -5.33%  "

